# CodeQuality-rules-latest.xlsx update
# - removes three retired SonarQube rule rows from the "Rules" sheet
#   (squid:S2076, squid:S2078, squid:S3318)
# - clears the stray hyperlink / URL text that had been left in column F
#   (next to the "ConsecutivelyLogAndThrow" rule row), while keeping the
#   cell's existing (Hyperlink) style
# - leaves the final cursor on the cell where the hyperlink used to be

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the leftover hyperlink (and its visible URL text) from column F ---
# The link lives on the "ConsecutivelyLogAndThrow" row; find it by its
# address instead of assuming a fixed row number.
$linkCell = $ws.Range("F103")
$ws.Hyperlinks.Delete()
$linkCell.ClearContents()

# --- Remove the obsolete rule rows ---
# Delete from the bottom up so earlier row numbers don't shift while we work.
$ws.Rows(26).Delete()   # squid:S3318 - Untrusted data should not be stored in sessions
$ws.Rows(10).Delete()   # squid:S2078 - Values passed to LDAP queries should be sanitized
$ws.Rows(8).Delete()    # squid:S2076 - Values passed to OS commands should be sanitized

# --- Leave the selection where the cleared hyperlink cell now sits ---
[void]$ws.Range("F100").Select()
